$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'30.100.10"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -2.00%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'1.829.33"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -3.26%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "'  +0.13%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'230.73"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -2.88%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 5).Value = "'  +0.11%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "'  -4.01%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).Value = "'  -6.88%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).Value = "'0.06260"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -4.61%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'1.830.52"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -2.64%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).Value = "'0.07387"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -0.82%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "'  -3.74%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'4.880"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -4.30%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'83.12"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -5.58%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'0.6140"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -7.60%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'30.076.66"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -1.90%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 5).Value = "'  +0.06%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'225.33"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -2.23%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'0.000007251"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -4.81%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 5).Value = "'  -6.50%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'0.9997"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +0.00%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'2.067.86"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -0.66%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'4.836"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -8.35%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).Value = "'5.815"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -6.35%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'9.142"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -3.07%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).Value = "'164.74"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -1.94%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'17.52"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -6.38%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'1.856"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -5.13%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).Value = "'0.1005"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -1.71%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 5).Value = "'  -1.58%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).Value = "'4.046"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -6.70%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = "'3.768"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -6.46%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 4).Value = "'0.04763"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  -5.78%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).Value = "'1.124"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -6.86%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 4).Value = "'0.7057"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -6.43%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 4).Value = "'2.699"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -0.41%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'0.01803"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -4.98%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'2.609"
$ws.Cells.Item(38, 4).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'0.8900"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -3.18%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'1.920"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -6.87%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 5).Value = "'  -0.17%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).Value = "'103.50"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -3.41%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'5.467"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -3.29%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'0.3984"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -6.97%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'6.930"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -6.92%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).Value = "'  -6.67%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'59.50"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -8.28%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'8.419"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -6.58%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).Value = "'32.64"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -4.23%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'0.05512"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -2.73%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).Value = "'1.365"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -7.83%  "
$ws.Cells.Item(51, 5).Style = "Normal"

